$d = $word.ActiveDocument

# The paragraph currently reads "Version 1." (0-based char offsets):
#   V  e  r  s  i  o  n     1  .
#   0  1  2  3  4  5  6  7  8  9
#
# Target reads "Version 2." but with "Version" split into two runs
# ("Versi"/"on") and the trailing "." moved to its own run positioned
# after the _GoBack bookmark.

# 1) Split "Version" into "Versi" | "on" without altering the visible
#    text or leaving stray formatting: temporarily bookmark the
#    mid-word split point, then delete that bookmark. Word leaves the
#    two text runs it created behind, un-merged, with no extra rPr.
$splitPoint = $d.Range(5, 5)
$d.Bookmarks.Add("__tmp_split", $splitPoint)
$d.Bookmarks.Item("__tmp_split").Delete()

# 2) "1" -> "2"
$digit = $d.Range(8, 9)
$digit.Text = "2"

# 3) Move the trailing "." so it sits after the _GoBack bookmark
#    instead of before it.
$period = $d.Range(9, 10)
$period.Text = ""
$paraEnd = $d.Paragraphs.Item(1).Range.End - 1
$insertPoint = $d.Range($paraEnd, $paraEnd)
$insertPoint.InsertAfter(".")
